$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that held George_Murray / Geoffrey (rows 4-5), which
# shifts Dick_Sheppard / Claus_Westermann up into rows 4-5.
$ws.Rows("4:5").Delete()

# Remove the trailing rows (Charles_Januarius_Acton, Carlo_Barberini,
# Cardinal_de_Bouillon) which are now at rows 6-8 after the shift above.
$ws.Rows("6:8").Delete()
